$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.925.66"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.647.84"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").Value = "'215.67"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.2576"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.06423"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "'19.72"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "'0.07776"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "'4.314"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "1.652.33"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "'0.5482"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "0.0₅7894"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'64.98"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "26.007.18"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "'198.11"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'6.067"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").Value = "'140.58"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "'0.1149"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "'6.902"
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "'1.241"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'0.05016"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'3.282"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'3.208"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'1.547"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'2.368"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'0.8965"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "'2.586"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.5537"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.131.50"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'5.661"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "'0.8159"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "'99.86"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("E44").Value = "  +8.25%  "
$ws.Range("D45").Value = "1.783.82"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'0.4543"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'55.36"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "'0.05098"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'0.09570"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("E51").Value = "  +0.14%  "
